# Apply updated crypto price/volume data per commit:
# "Updated cryptos list on Wed Jan 10 22:59:55 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '46.885.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.546.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +9.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.608'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.570'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +12.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0829'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.94'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +12.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.931.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.01%  '
$ws.Range('E14').Value = '  +2.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.571.34'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +10.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.899'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +12.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.04'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +10.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '46.787.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +12.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0994'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.84%  '
$ws.Range('E24').Value = '  +6.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.15%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '42.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +19.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.50'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.39%  '
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.01'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.94'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0844'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +24.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '150.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.122'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.83%  '
$ws.Range('E38').Value = '  +4.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.29'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0326'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +13.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.003.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.996'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '94.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.59'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +37.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.202'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.59%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.93%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.86'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +12.22%  '
